$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.399.07"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.847.17"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.91"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6289"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.46"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").Value = "1.850.92"
$ws.Range("E12").Value = "  -6.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.002"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001089"
$ws.Range("E14").Value = "  +9.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6789"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.71"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "2.098.45"
$ws.Range("E17").Value = "  -7.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.171"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "29.421.00"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.98"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.441"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.46"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1397"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.364"
$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.303"
$ws.Range("E30").Value = "  +4.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05597"
$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.034"
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7095"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "1.231.98"
$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.769"
$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.454"
$ws.Range("E41").Value = "  +5.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9068"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.85"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.04"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000121"
$ws.Range("E46").Value = "  +2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.174"
$ws.Range("E47").Value = "  +1.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4019"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.957"
$ws.Range("E49").Value = "  -2.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.678"
$ws.Range("E50").Value = "  -1.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1122"
$ws.Range("E51").Value = "  -0.44%  "
